$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.643.33"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "3.048.34"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'213.21"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "'608.85"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("D7").Value = "'1.06"
$ws.Range("E7").Value = "  +18.73%  "
$ws.Range("D8").Value = "'0.344"
$ws.Range("E8").Value = "  -10.05%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "3.049.93"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("D11").Value = "'0.705"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "'0.193"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "'0.0000235"
$ws.Range("E13").Value = "  -7.87%  "
$ws.Range("D14").Value = "'5.45"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'33.68"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "89.501.38"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "3.635.53"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "3.058.70"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'3.62"
$ws.Range("E19").Value = "  -4.01%  "
$ws.Range("D20").Value = "'13.92"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'444.72"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'0.0000198"
$ws.Range("E22").Value = "  -15.70%  "
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").Value = "'5.38"
$ws.Range("E23").Value = "  +4.26%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'8.69"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'5.66"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "'90.41"
$ws.Range("E26").Value = "  +7.08%  "
$ws.Range("D27").Value = "'11.68"
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("D28").Value = "3.259.47"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'9.04"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.156"
$ws.Range("E31").Value = "  -8.11%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'27.78"
$ws.Range("E32").Value = "  +19.44%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.199"
$ws.Range("E33").Value = "  +39.96%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'0.928"
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'491.13"
$ws.Range("E35").Value = "  -5.86%  "
$ws.Range("D36").Value = "'0.143"
$ws.Range("E36").Value = "  +3.28%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "'1.86"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'6.71"
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("D39").Value = "'1.26"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  -12.70%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.417"
$ws.Range("E41").Value = "  +11.08%  "
$ws.Range("D42").Value = "'22.14"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0839"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "'1.90"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "  +9.20%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'147.41"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'4.54"
$ws.Range("E48").Value = "  +8.62%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.674"
$ws.Range("E49").Value = "  +10.02%  "
$ws.Range("D50").Value = "'44.57"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.08%  "
